# Applies the values for columns H (F2_TieBreak_Nuevos_Clientes) and
# I (F3_Pedidos_Por_Dia) for rows 2-13, and updates the active selection
# on the worksheet from G14 to I14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hValues = @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$iValues = @(100, 100, 99, 99, 99, 99, 99, 99, 99, 99, 99, 99)

for ($i = 0; $i -lt $hValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
}

$ws.Range("I14").Select()
